$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '42.773.80'
$ws.Cells.Item(2, 5).Value = '  -1.78%  '
$ws.Cells.Item(3, 4).Value = '2.566.76'
$ws.Cells.Item(3, 5).Value = '  -1.12%  '
$ws.Cells.Item(4, 5).Value = '  +0.18%  '
$ws.Cells.Item(5, 4).Value = '302.26'
$ws.Cells.Item(5, 5).Value = '  +0.54%  '
$ws.Cells.Item(6, 4).Value = '92.86'
$ws.Cells.Item(6, 5).Value = '  -3.58%  '
$ws.Cells.Item(7, 5).Value = '  -0.70%  '
$ws.Cells.Item(8, 5).Value = '  +0.01%  '
$ws.Cells.Item(9, 5).Value = '  -2.38%  '
$ws.Cells.Item(10, 4).Value = '35.98'
$ws.Cells.Item(10, 5).Value = '  -2.23%  '
$ws.Cells.Item(11, 5).Value = '  -0.68%  '
$ws.Cells.Item(12, 4).Value = '7.65'
$ws.Cells.Item(12, 5).Value = '  -2.15%  '
$ws.Cells.Item(13, 5).Value = '  +6.67%  '
$ws.Cells.Item(14, 4).Value = '2.554.26'
$ws.Cells.Item(14, 5).Value = '  -1.20%  '
$ws.Cells.Item(15, 4).Value = '0.882'
$ws.Cells.Item(15, 5).Value = '  -1.16%  '
$ws.Cells.Item(16, 5).Value = '  -0.91%  '
$ws.Cells.Item(17, 4).Value = '42.863.93'
$ws.Cells.Item(17, 5).Value = '  -1.44%  '
$ws.Cells.Item(18, 4).Value = '0.0₃0988'
$ws.Cells.Item(18, 5).Value = '  +1.08%  '
$ws.Cells.Item(19, 4).Value = '12.75'
$ws.Cells.Item(19, 5).Value = '  +3.34%  '
$ws.Cells.Item(20, 4).Value = '6.62'
$ws.Cells.Item(20, 5).Value = '  -0.80%  '
$ws.Cells.Item(21, 4).Value = '71.79'
$ws.Cells.Item(21, 5).Value = '  -1.47%  '
$ws.Cells.Item(22, 4).Value = '253.12'
$ws.Cells.Item(22, 5).Value = '  -5.06%  '
$ws.Cells.Item(23, 4).Value = '2.95'
$ws.Cells.Item(24, 5).Value = '  -4.36%  '
$ws.Cells.Item(25, 4).Value = '28.78'
$ws.Cells.Item(25, 5).Value = '  -2.16%  '
$ws.Cells.Item(26, 5).Value = '  -0.14%  '
$ws.Cells.Item(27, 4).Value = '10.30'
$ws.Cells.Item(27, 5).Value = '  +0.29%  '
$ws.Cells.Item(28, 4).Value = '36.92'
$ws.Cells.Item(28, 5).Value = '  -1.65%  '
$ws.Cells.Item(29, 5).Value = '  -4.20%  '
$ws.Cells.Item(30, 5).Value = '  -1.02%  '
$ws.Cells.Item(31, 5).Value = '  +1.65%  '
$ws.Cells.Item(32, 5).Value = '  -3.95%  '
$ws.Cells.Item(33, 5).Value = '  -6.24%  '
$ws.Cells.Item(34, 5).Value = '  -1.20%  '
$ws.Cells.Item(35, 4).Value = '0.0798'
$ws.Cells.Item(35, 5).Value = '  -2.15%  '
$ws.Cells.Item(36, 4).Value = '18.36'
$ws.Cells.Item(36, 5).Value = '  +9.44%  '
$ws.Cells.Item(37, 5).Value = '  -3.88%  '
$ws.Cells.Item(38, 5).Value = '  -1.32%  '
$ws.Cells.Item(39, 4).Value = '23.44'
$ws.Cells.Item(39, 5).Value = '  -4.79%  '
$ws.Cells.Item(40, 5).Value = '  +33.17%  '
$ws.Cells.Item(41, 5).Value = '  -1.36%  '
$ws.Cells.Item(42, 4).Value = '3.40'
$ws.Cells.Item(42, 5).Value = '  -4.40%  '
$ws.Cells.Item(43, 4).Value = '3.87'
$ws.Cells.Item(43, 5).Value = '  +0.52%  '
$ws.Cells.Item(44, 4).Value = '2.083.59'
$ws.Cells.Item(44, 5).Value = '  +1.93%  '
$ws.Cells.Item(46, 4).Value = '9.25'
$ws.Cells.Item(46, 5).Value = '  +2.19%  '
$ws.Cells.Item(47, 4).Value = '84.84'
$ws.Cells.Item(47, 5).Value = '  -4.44%  '
$ws.Cells.Item(48, 2).Value = 'ordi'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Cells.Item(48, 4).Value = '75.96'
$ws.Cells.Item(48, 5).Value = '  +9.64%  '
$ws.Cells.Item(49, 2).Value = 'Aave'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(49, 4).Value = '106.54'
$ws.Cells.Item(49, 5).Value = '  +0.52%  '
$ws.Cells.Item(50, 4).Value = '2.819.40'
$ws.Cells.Item(50, 5).Value = '  -0.64%  '
$ws.Cells.Item(51, 2).Value = 'Algorand'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(51, 4).Value = '0.191'
$ws.Cells.Item(51, 5).Value = '  +0.06%  '
